# "actualizacao da sessao 8 chimanimani"
# Update the attendance/engagement record for producer C_1074 (row 2) on
# Sheet1: presenca (M), nrSessao (N) and apresentacao_photovoice (P) get
# new values for this session, and nivel_engajamento (Q) is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "NAO"   # presenca: SIM -> NAO
$ws.Range("N2").Value = 6       # nrSessao: 2 -> 6
$ws.Range("P2").Value = "NAO"   # apresentacao_photovoice: (blank) -> NAO
$ws.Range("Q2").Value = ""      # nivel_engajamento: ENGAJAD@ -> (blank)
